# Add a "Total:" label and a jxls SUM formula placeholder below the
# detail row, matching the sample-jxls1.xlsx template update.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 14: D14 = "Total:", E14 = "$[SUM(E11)]"
# (write E14's string first so the shared-string table gets the same
#  ordering as the authored workbook: $[SUM(E11)] before Total:)
$ws.Range("E14").Value = "`$[SUM(E11)]"
$ws.Range("D14").Value = "Total:"

# Match the header row's style (bold font) for the new cells.
$ws.Range("D14:E14").Font.Bold = $true

# Update selection to reflect the newly active range.
$ws.Range("D14:E14").Select()
